# This script inserts two new daily price records (rows) into the
# "Plátano" consolidated sheet, at position 664/665, pushing every
# subsequent row down by two. The two new rows carry a brand-new
# observation date (2021-11-22 / serial 44522) for variety
# "Sin especificar" with qualities "Pintón" and "Primera Pintón".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 664, shifting the existing data (old rows
# 664..776) down to 666..778.
$ws.Rows("664:665").Insert()

# --- New row 664 ---------------------------------------------------
$ws.Cells.Item(664, 1).Value = 6
$ws.Cells.Item(664, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(664, 3).Value = "Metropolitana"
$ws.Cells.Item(664, 4).Value = 44522
$ws.Cells.Item(664, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(664, 5).Value = 13
$ws.Cells.Item(664, 6).Value = "Fruta"
$ws.Cells.Item(664, 7).Value = 100108
$ws.Cells.Item(664, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(664, 9).Value = 100108006
$ws.Cells.Item(664, 10).Value = "Plátano"
$ws.Cells.Item(664, 11).Value = "Sin especificar"
$ws.Cells.Item(664, 12).Value = "Pintón"
$ws.Cells.Item(664, 13).Value = 760
$ws.Cells.Item(664, 14).Value = 19000
$ws.Cells.Item(664, 15).Value = 20000
$ws.Cells.Item(664, 16).Value = 19737
$ws.Cells.Item(664, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(664, 18).Value = "Ecuador"
$ws.Cells.Item(664, 19).Value = 987
$ws.Cells.Item(664, 20).Value = 20

# --- New row 665 ---------------------------------------------------
$ws.Cells.Item(665, 1).Value = 6
$ws.Cells.Item(665, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(665, 3).Value = "Metropolitana"
$ws.Cells.Item(665, 4).Value = 44522
$ws.Cells.Item(665, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(665, 5).Value = 13
$ws.Cells.Item(665, 6).Value = "Fruta"
$ws.Cells.Item(665, 7).Value = 100108
$ws.Cells.Item(665, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(665, 9).Value = 100108006
$ws.Cells.Item(665, 10).Value = "Plátano"
$ws.Cells.Item(665, 11).Value = "Sin especificar"
$ws.Cells.Item(665, 12).Value = "Primera Pintón"
$ws.Cells.Item(665, 13).Value = 1000
$ws.Cells.Item(665, 14).Value = 21000
$ws.Cells.Item(665, 15).Value = 22000
$ws.Cells.Item(665, 16).Value = 21400
$ws.Cells.Item(665, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(665, 18).Value = "Ecuador"
$ws.Cells.Item(665, 19).Value = 1070
$ws.Cells.Item(665, 20).Value = 20
